$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new journal entry row (row 16), matching the style of the
# previous entry rows (date format on A, number on B, wrapped text on C).

# A16: date 22.2.2022 (serial 44614), same style as A15 (date format)
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A16").Value = 44614

# B16: formula 60+120 -> 180 minutes, same style as B15
$ws.Range("B15").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("B16").Formula = "=60+120"

# C16: new entry text, same style as C15 (wrapped text)
$ws.Range("C15").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C16").Value = "Nba komponentin proton jatkamista loppuun"

# Recalculate so the total-hours array formula in F1 picks up the new row
$excel.CalculateFullRebuild() | Out-Null

# Leave the selection where it ended up after entering the new row
$ws.Range("B17").Select() | Out-Null
